# feat: add 2022-Q1 data
#
# - Insert a new "2022-Q1" worksheet (fund-holding detail, same layout as the
#   existing 2021-Q2 / 2021-Q3 / 2021-Q4 sheets) right before the "总计" sheet.
# - Prepend a matching summary row to the "总计" sheet.

$wb = $excel.ActiveWorkbook
$sheets = $wb.Worksheets

$totalSheetName = "总计"
$totalSheet = $sheets.Item($totalSheetName)
$formatSrc = $sheets.Item("2021-Q4")

# ---------------------------------------------------------------------------
# 1. New "2022-Q1" detail sheet, inserted right before "总计"
# ---------------------------------------------------------------------------
$q1 = $sheets.Add($totalSheet)
$q1.Name = "2022-Q1"

# Bring over the header/row formatting (fonts, borders, alignment) used by
# the other quarterly sheets so the new sheet matches their look.
$formatSrc.Range("A1:H2").Copy()
$q1.Range("A1:H2").PasteSpecial(-4122)
$formatSrc.Range("A2:H2").Copy()
$q1.Range("A3:H5").PasteSpecial(-4122)

# Header row
$q1.Range("B1").Value = "基金代码"
$q1.Range("C1").Value = "基金名称"
$q1.Range("D1").Value = "基金规模"
$q1.Range("E1").Value = "股票总仓位"
$q1.Range("F1").Value = "仓位占比"
$q1.Range("G1").Value = "持有市值(亿元)"
$q1.Range("H1").Value = "仓位排名"

# Fund holding rows for 2022-Q1
$q1Data = @(
    @{ Idx=0; Code="206009"; Name="鹏华新兴产业混合";     Size="44.95"; Pos="90.17"; Pct="8.20"; Mv="3.6859"; Rank=1 },
    @{ Idx=1; Code="011471"; Name="鹏华致远成长混合A";    Size="2.19";  Pos="61.03"; Pct="2.62"; Mv="0.0574"; Rank=9 },
    @{ Idx=2; Code="000166"; Name="中海信息产业精选混合"; Size="0.73";  Pos="91.91"; Pct="3.86"; Mv="0.0282"; Rank=7 },
    @{ Idx=3; Code="011472"; Name="鹏华致远成长混合C";    Size="0.08";  Pos="61.03"; Pct="2.62"; Mv="0.0021"; Rank=9 }
)

$row = 2
foreach ($item in $q1Data) {
    $q1.Range("A$row").Value = $item.Idx

    $codeCell = $q1.Range("B$row")
    $codeCell.NumberFormat = "@"
    $codeCell.Value = $item.Code

    $q1.Range("C$row").Value = $item.Name

    $sizeCell = $q1.Range("D$row")
    $sizeCell.NumberFormat = "@"
    $sizeCell.Value = $item.Size

    $posCell = $q1.Range("E$row")
    $posCell.NumberFormat = "@"
    $posCell.Value = $item.Pos

    $pctCell = $q1.Range("F$row")
    $pctCell.NumberFormat = "@"
    $pctCell.Value = $item.Pct

    $mvCell = $q1.Range("G$row")
    $mvCell.NumberFormat = "@"
    $mvCell.Value = $item.Mv

    $q1.Range("H$row").Value = $item.Rank

    $row = $row + 1
}

# ---------------------------------------------------------------------------
# 2. Prepend the 2022-Q1 summary row to "总计"
# ---------------------------------------------------------------------------
# Worksheet handles are positional, and adding $q1 above shifted everything
# at/after its insertion point, so re-resolve "总计" by name before touching it.
$totalSheet = $sheets.Item($totalSheetName)
$totalSheet.Rows("2:2").Insert(-4121)

# Restore column-A styling (index style) on the newly inserted row, then
# clear the stray formatting Insert() carried over into B:D.
$totalSheet.Range("A3").Copy()
$totalSheet.Range("A2").PasteSpecial(-4122)
$totalSheet.Range("B2:D2").ClearFormats()

$totalSheet.Range("A2").Value = 0
$totalSheet.Range("B2").Value = "2022-Q1"
$totalSheet.Range("C2").Value = 4
$totalSheet.Range("D2").Value = 3.77

# The A column is a 0-based row index, not literal data that should carry
# along with the shifted rows -- renumber it to stay sequential.
$totalSheet.Range("A3").Value = 1
$totalSheet.Range("A4").Value = 2
$totalSheet.Range("A5").Value = 3
